# "#deleted additional unused variables and removed corresponding SSC INPUTs"
#
# Appends 8 more "Deleted variable" rows to the bottom of the
# "SAM Variable Changes" sheet (rows 49-56), documenting SSC inputs that
# were removed from the (now-deleted) "Direct Steam Tower Receiver" page:
#   rec_htf_flow, ref_htf_flow, P_cold_tank, P_tower_conv, P_tower_rad,
#   P_htf_pump, night_recirc, mode
# Each new row follows the same shape/style as the existing rows directly
# above it (rows 44-48, also "Deleted variable" / "not used" rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SAM Variable Changes")

# Old (deleted) variable names, one per new row, column C
$deletedNames = @(
    "rec_htf_flow",
    "ref_htf_flow",
    "P_cold_tank",
    "P_tower_conv",
    "P_tower_rad",
    "P_htf_pump",
    "night_recirc",
    "mode"
)

$startRow = 49
$srcRow   = 48   # last pre-existing row of the same kind; used as a format template

# Clone the per-column formatting from row 48 down onto the new rows so the
# new cells pick up the same styles (fill/alignment) used by the rest of
# this "deleted variable" block.
foreach ($col in @(1, 2, 3, 5, 6, 7, 8)) {
    $ws.Cells.Item($srcRow, $col).Copy() | Out-Null
    for ($i = 0; $i -lt $deletedNames.Length; $i++) {
        $row = $startRow + $i
        $ws.Cells.Item($row, $col).PasteSpecial(-4122) | Out-Null   # xlPasteFormats
    }
}
$excel.CutCopyMode = $false

for ($i = 0; $i -lt $deletedNames.Length; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = "Deleted variable"             # A: Type
    $ws.Cells.Item($row, 2).Value = "number"                       # B: Variable Type
    $ws.Cells.Item($row, 3).Value = $deletedNames[$i]               # C: Old Name
    $ws.Cells.Item($row, 5).Value = "Direct Steam Tower Receiver"  # E: Input Page
    $ws.Cells.Item($row, 6).Value = "not used"                     # F: Default Value (if new) or reason deleted
    $ws.Cells.Item($row, 7).Value = "N/A"                          # G: Handled in Version Upgrader?
    $ws.Cells.Item($row, 8).Value = "Ty"                           # H: (initials)
}

# Reflect the post-edit selection/scroll state (user clicked just past the
# last newly-added row).
$ws.Range("A57").Select()
